$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the existing data block (rows 4-7) twice, into rows 8-11 and 12-15.
# This carries over cell styles (bold names, hyperlink-style email column,
# date/number formats) exactly like the source rows.
$ws.Range("A4:L7").Copy($ws.Range("A8:L11"))
$ws.Range("A4:L7").Copy($ws.Range("A12:L15"))

# Tweak the EMAIL_ADDRESS column for the first duplicated block ("...1@...")
$ws.Range("C11").Value = "dsjhds1@fdkfjdk.com"
$ws.Range("C10").Value = "luka.doncic1@me.com"
$ws.Range("C9").Value = "ian.intermediate1@gmail.com"
$ws.Range("C8").Value = "nick.newbie1@gmail.com"

# Tweak the EMAIL_ADDRESS column for the second duplicated block ("...2@...")
$ws.Range("C15").Value = "dsjhds2@fdkfjdk.com"
$ws.Range("C14").Value = "luka.doncic2@me.com"
$ws.Range("C13").Value = "ian.intermediate2@gmail.com"
$ws.Range("C12").Value = "nick.newbie2@gmail.com"

# Mark the two blocks with a short tag in column A ("a" / "b")
$ws.Range("A8").Value = "a"
$ws.Range("A9").Value = "a"
$ws.Range("A10").Value = "a"
$ws.Range("A11").Value = "a"
$ws.Range("A12").Value = "b"
$ws.Range("A13").Value = "b"
$ws.Range("A14").Value = "b"
$ws.Range("A15").Value = "b"

# Re-create the mailto hyperlinks for the new EMAIL_ADDRESS cells, in row
# order (C8..C15), so rId6..rId13 line up with the respective rows.
$ws.Hyperlinks.Add($ws.Range("C8"), "mailto:nick.newbie1@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C9"), "mailto:ian.intermediate1@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C10"), "mailto:luka.doncic1@me.com")
$ws.Hyperlinks.Add($ws.Range("C11"), "mailto:dsjhds1@fdkfjdk.com")
$ws.Hyperlinks.Add($ws.Range("C12"), "mailto:nick.newbie2@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C13"), "mailto:ian.intermediate2@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C14"), "mailto:luka.doncic2@me.com")
$ws.Hyperlinks.Add($ws.Range("C15"), "mailto:dsjhds2@fdkfjdk.com")

# Update the saved selection/active cell, matching the final cursor position.
$null = $ws.Range("H20").Select()
